$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 60, pushing rows 60-171 down to 61-172.
$ws.Rows("60:60").Insert()

# Populate the newly inserted row 60 with its data.
$ws.Range("A60").Value = 11
$ws.Range("B60").Value = "Vega Monumental Concepción"
$ws.Range("C60").Value = "Bíobío"
$ws.Range("D60").Value = 44498
$ws.Range("E60").Value = 8
$ws.Range("F60").Value = 100112009
$ws.Range("G60").Value = "Acelga"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 450
$ws.Range("K60").Value = 600
$ws.Range("L60").Value = 650
$ws.Range("M60").Value = 622
$ws.Range("N60").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O60").Value = "Región del Maule"
$ws.Range("P60").Value = 415
$ws.Range("Q60").Value = 1.5
$ws.Range("R60").Value = "Hortaliza"
